# Auto-generated edit script: updates cryptos list (prices/volumes) and
# fixes two pairs of rows whose coin order flipped (29/30 and 42/43/44/45),
# plus row 51 Monero -> Arweave, matching the "Updated cryptos list" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "64.283.12"
$ws.Range("E2").Value = "  +1.86%  "

# Row 3
$ws.Range("D3").Value = "3.311.66"
$ws.Range("E3").Value = "  +6.00%  "

# Row 4
$ws.Range("E4").Value = "  +0.11%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "600.67"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.39%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.13"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +5.98%  "

# Row 7
$ws.Range("E7").Value = "  +0.04%  "

# Row 8
$ws.Range("D8").Value = "3.304.39"
$ws.Range("E8").Value = "  +6.01%  "

# Row 9
$ws.Range("E9").Value = "  +1.01%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.150"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +2.56%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.49"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +3.59%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000250"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +0.41%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.18"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +3.45%  "

# Row 15
$ws.Range("D15").Value = "3.855.86"
$ws.Range("E15").Value = "  +5.87%  "

# Row 16
$ws.Range("E16").Value = "  +1.49%  "

# Row 17
$ws.Range("D17").Value = "3.310.61"
$ws.Range("E17").Value = "  +6.41%  "

# Row 18
$ws.Range("D18").Value = "64.356.07"
$ws.Range("E18").Value = "  +2.05%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.93"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +3.30%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "487.36"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +2.41%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.41"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +1.75%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.747"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +7.15%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.09"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +6.01%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.65"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +4.49%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "85.16"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -2.83%  "

# Row 26
$ws.Range("E26").Value = "  -0.17%  "

# Row 27
$ws.Range("E27").Value = "  +3.25%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.35"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +3.95%  "

# Row 29
$ws.Range("B29").Value = "FirstDigitalUSD"
$ws.Range("C29").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +0.10%  "

# Row 30
$ws.Range("B30").Value = "NEARProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.26"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +1.01%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.18"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +5.37%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "28.48"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +4.37%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.108"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +0.37%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.59"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +2.06%  "

# Row 35
$ws.Range("E35").Value = "  +2.53%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.05"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +3.42%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "53.33"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +2.66%  "

# Row 38
$ws.Range("D38").Value = "0.0₃0745"
$ws.Range("E38").Value = "  +4.57%  "

# Row 39
$ws.Range("E39").Value = "  +3.14%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "431.79"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +2.45%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.82"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +4.28%  "

# Row 42
$ws.Range("B42").Value = "Cosmos"
$ws.Range("C42").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.50"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +2.83%  "

# Row 43
$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").Value = "3.020.55"
$ws.Range("E43").Value = "  +5.63%  "

# Row 44
$ws.Range("B44").Value = "Kaspa"
$ws.Range("C44").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.111"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -5.40%  "

# Row 45
$ws.Range("B45").Value = "TheGraph"
$ws.Range("C45").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.273"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +5.96%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.26"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +7.16%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "26.46"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +4.24%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.35"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +2.56%  "

# Row 49
$ws.Range("E49").Value = "  +0.11%  "

# Row 50
$ws.Range("E50").Value = "  +1.79%  "

# Row 51
$ws.Range("B51").Value = "Arweave"
$ws.Range("C51").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "34.93"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +13.82%  "

